# Auto-generated script: apply scheduled-runner market data updates
# to the "Coeurl" profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 8423.727999999999
$ws.Range("I9").Value = 11519.875
$ws.Range("J9").Value = 167.33333
$ws.Range("K9").Value = 11519.875
$ws.Range("L9").Value = 167.33333
$ws.Range("M9").Value = -11350.875
$ws.Range("N9").Value = -505.33333

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 714322.5600000001
$ws.Range("J17").Value = 714322.5600000001
$ws.Range("L17").Value = 2142967.68
$ws.Range("N17").Value = -2143303.68

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 4529.5557
$ws.Range("I86").Value = 3846.375
$ws.Range("K86").Value = 3846.375
$ws.Range("M86").Value = -2723.375

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 4529.5557
$ws.Range("I89").Value = 3846.375
$ws.Range("K89").Value = 19231.875
$ws.Range("M89").Value = -13615.875

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 54321.316
$ws.Range("J112").Value = 64240.375
$ws.Range("L112").Value = 192721.125
$ws.Range("N112").Value = -194937.125

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H133").Value = 76487.5
$ws.Range("J133").Value = 76487.5
$ws.Range("L133").Value = 76487.5
$ws.Range("N133").Value = -86607.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").Value = ""

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 2036.1111
$ws.Range("I141").Value = 1861.7646
$ws.Range("K141").Value = 5585.293799999999
$ws.Range("M141").Value = -405.2937999999995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 91006440
$ws.Range("I2").Value = 125131290
$ws.Range("J2").Value = 6831
$ws.Range("K2").Value = 125131290
$ws.Range("L2").Value = 6831
$ws.Range("M2").Value = -125131177
$ws.Range("N2").Value = -7057

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 22413.6
$ws.Range("I45").Value = 26274.75
$ws.Range("J45").Value = 6969
$ws.Range("K45").Value = 26274.75
$ws.Range("L45").Value = 6969
$ws.Range("M45").Value = -25897.75
$ws.Range("N45").Value = -7723

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2793.7014
$ws.Range("I61").Value = 2503.585
$ws.Range("J61").Value = 3892
$ws.Range("K61").Value = 2503.585
$ws.Range("L61").Value = 3892
$ws.Range("M61").Value = -2291.585
$ws.Range("N61").Value = -4316

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 7345.619
$ws.Range("I110").Value = 8217.267
$ws.Range("J110").Value = 5166.5
$ws.Range("K110").Value = 8217.267
$ws.Range("L110").Value = 5166.5
$ws.Range("M110").Value = -6172.267
$ws.Range("N110").Value = -9256.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 91006440
$ws.Range("I116").Value = 125131290
$ws.Range("J116").Value = 6831
$ws.Range("K116").Value = 125131290
$ws.Range("L116").Value = 6831
$ws.Range("M116").Value = -125128996
$ws.Range("N116").Value = -11419

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 3392.0908
$ws.Range("I122").Value = 3349.3157
$ws.Range("J122").Value = 3663
$ws.Range("K122").Value = 10047.9471
$ws.Range("L122").Value = 10989
$ws.Range("M122").Value = -7597.947100000001
$ws.Range("N122").Value = -15889

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2585.3635
$ws.Range("I132").Value = 2522.375
$ws.Range("J132").Value = 2753.3333
$ws.Range("K132").Value = 7567.125
$ws.Range("L132").Value = 8259.999899999999
$ws.Range("M132").Value = -5037.125
$ws.Range("N132").Value = -13319.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2793.7014
$ws.Range("I136").Value = 2503.585
$ws.Range("J136").Value = 3892
$ws.Range("K136").Value = 7510.755
$ws.Range("L136").Value = 11676
$ws.Range("M136").Value = -4960.755
$ws.Range("N136").Value = -16776

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 91006440
$ws.Range("I3").Value = 125131290
$ws.Range("J3").Value = 6831
$ws.Range("K3").Value = 125131290
$ws.Range("L3").Value = 6831
$ws.Range("M3").Value = -125131176
$ws.Range("N3").Value = -7059

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3041.7878
$ws.Range("I20").Value = 2750.1667
$ws.Range("J20").Value = 3208.4285
$ws.Range("K20").Value = 2750.1667
$ws.Range("L20").Value = 3208.4285
$ws.Range("M20").Value = -2503.1667
$ws.Range("N20").Value = -3702.4285

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 864.75
$ws.Range("I94").Value = 1326.6666
$ws.Range("J94").Value = 402.83334
$ws.Range("K94").Value = 1326.6666
$ws.Range("L94").Value = 402.83334
$ws.Range("M94").Value = -875.6666
$ws.Range("N94").Value = -1304.83334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2808.8333
$ws.Range("I107").Value = 2388
$ws.Range("J107").Value = 3109.4285
$ws.Range("K107").Value = 2388
$ws.Range("L107").Value = 3109.4285
$ws.Range("M107").Value = -468
$ws.Range("N107").Value = -6949.4285

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").Value = ""

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1433.6666
$ws.Range("I105").Value = 1433.6666
$ws.Range("K105").Value = 1433.6666
$ws.Range("M105").Value = 313.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 13893.289
$ws.Range("I134").Value = 7922.0645
$ws.Range("J134").Value = 40337.285
$ws.Range("K134").Value = 23766.1935
$ws.Range("L134").Value = 121011.855
$ws.Range("M134").Value = -21231.1935
$ws.Range("N134").Value = -126081.855

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 59.61905
$ws.Range("J2").Value = 35.333332
$ws.Range("L2").Value = 211.999992
$ws.Range("N2").Value = -437.999992

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 19998.637
$ws.Range("I56").Value = 19998.637
$ws.Range("K56").Value = 19998.637
$ws.Range("M56").Value = -19468.637

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 3697.8
$ws.Range("J70").Value = 3872.25
$ws.Range("L70").Value = 11616.75
$ws.Range("N70").Value = -12246.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H73").Value = 3697.8
$ws.Range("J73").Value = 3872.25
$ws.Range("L73").Value = 11616.75
$ws.Range("N73").Value = -13800.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 13046.066
$ws.Range("I87").Value = 8807.583000000001
$ws.Range("K87").Value = 26422.749
$ws.Range("M87").Value = -25174.749

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H88").Value = 33334394
$ws.Range("J88").Value = 33334394
$ws.Range("L88").Value = 100003182
$ws.Range("N88").Value = -100004038

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H90").Value = 13046.066
$ws.Range("I90").Value = 8807.583000000001
$ws.Range("K90").Value = 79268.247
$ws.Range("M90").Value = -73028.247

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H91").Value = 33334394
$ws.Range("J91").Value = 33334394
$ws.Range("L91").Value = 100003182
$ws.Range("N91").Value = -100006146

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H115").Value = 13275.333
$ws.Range("I115").Value = 13275.333
$ws.Range("K115").Value = 39825.999
$ws.Range("M115").Value = -38650.999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H59").Value = 15000
$ws.Range("J59").Value = 15000
$ws.Range("L59").Value = 15000
$ws.Range("N59").Value = -16166

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").Value = ""

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").Value = ""

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").Value = ""

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1766.762
$ws.Range("I122").Value = 1169.4
$ws.Range("K122").Value = 3508.2
$ws.Range("M122").Value = -1058.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H138").Value = 59429
$ws.Range("J138").Value = 59429
$ws.Range("L138").Value = 59429
$ws.Range("N138").Value = -69709

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 12430
$ws.Range("I43").Value = 8795.556
$ws.Range("K43").Value = 8795.556
$ws.Range("M43").Value = -8602.556

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 5562.5557
$ws.Range("I122").Value = 4611.1763
$ws.Range("J122").Value = 7179.9
$ws.Range("K122").Value = 13833.5289
$ws.Range("L122").Value = 21539.7
$ws.Range("M122").Value = -11383.5289
$ws.Range("N122").Value = -26439.7

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 630.2857
$ws.Range("I107").Value = 667.63635
$ws.Range("J107").Value = 493.33334
$ws.Range("K107").Value = 2002.90905
$ws.Range("L107").Value = 1480.00002
$ws.Range("M107").Value = -82.90904999999998
$ws.Range("N107").Value = -5320.000019999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2102.125
$ws.Range("I122").Value = 2051.5454
$ws.Range("J122").Value = 2340.5715
$ws.Range("K122").Value = 6154.6362
$ws.Range("L122").Value = 7021.7145
$ws.Range("M122").Value = -3704.6362
$ws.Range("N122").Value = -11921.7145

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 5487.5
$ws.Range("I126").Value = 4983.3335
$ws.Range("K126").Value = 14950.0005
$ws.Range("M126").Value = -12480.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H137").Value = 100666.5
$ws.Range("J137").Value = 100666.5
$ws.Range("L137").Value = 100666.5
$ws.Range("N137").Value = -110866.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H140").Value = 74697.25
$ws.Range("J140").Value = 74697.25
$ws.Range("L140").Value = 74697.25
$ws.Range("N140").Value = -85057.25
